$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column E ("reviews_count") entirely, shifting columns F:K left by one
# (F->E, G->F, H->G, I->H, J->I, K->J).
$ws.Columns.Item(5).Delete()
